$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Small precision corrections to existing historical values ---
$ws.Range("B35").Value = 4.5468317896932238
$ws.Range("B120").Value = 1.7585324826803159
$ws.Range("B121").Value = 3.4167832824663269

# --- Fill in the previously-blank forecast rows 131-138 ---
$ws.Range("B131").Value = 4.9531529009810278
$ws.Range("C131").Value = 4
$ws.Range("B132").Value = 3.0489989112225038
$ws.Range("C132").Value = 3
$ws.Range("B133").Value = 2.7577358723358074
$ws.Range("C133").Value = 3
$ws.Range("B134").Value = 3.411146524749809
$ws.Range("C134").Value = 2.5
$ws.Range("B135").Value = 3.2004613492993212
$ws.Range("C135").Value = 4
$ws.Range("B136").Value = 4.2699463713570358
$ws.Range("C136").Value = 3.5
$ws.Range("B137").Value = 4.3415269084839139
$ws.Range("C137").Value = 3.5
$ws.Range("B138").Value = 2.1644386352046396
$ws.Range("C138").Value = 3

# --- Append new month-end rows 142-147, cloning the style of row 141 ---
$newDates = @(45536, 45566, 45597, 45627, 45658, 45689)
$row = 141
foreach ($d in $newDates) {
    $row = $row + 1
    $ws.Range("A141").Copy($ws.Range("A$row"))
    $ws.Range("A$row").Value = $d
    $ws.Range("B141").Copy($ws.Range("B$row"))
    $ws.Range("C141").Copy($ws.Range("C$row"))
}
